# Adds the "2022-Q3" quarterly holdings sheet and updates the "总计"
# (summary) sheet so the new quarter appears at the top of the table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计" (so tab order is
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4).
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row (bold, thin border, centered, top-aligned - matches the style
# used for header rows / index column elsewhere in this workbook).
$q3Header = $q3.Range("B1:H1")
$q3Header.Font.Bold = $true
$q3Header.Borders.LineStyle = 1
$q3Header.HorizontalAlignment = -4108
$q3Header.VerticalAlignment = -4160

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Index column (A2:A11) uses the same look as the header row.
$q3Index = $q3.Range("A2:A11")
$q3Index.Font.Bold = $true
$q3Index.Borders.LineStyle = 1
$q3Index.HorizontalAlignment = -4108
$q3Index.VerticalAlignment = -4160

# Fund holdings data for 2022-Q3, ordered exactly as in the source sheet.
$q3rows = @(
    @("001403", "招商国企改革主题混合", "1.97", "85.29", "6.40", "0.1261", 3),
    @("217001", "招商安泰混合", "4.18", "75.08", "2.46", "0.1028", 9),
    @("000030", "长城核心优选灵活配置混合", "1.58", "93.95", "5.26", "0.0831", 9),
    @("519677", "银河定投宝腾讯济安指数", "2.88", "92.40", "2.48", "0.0714", 1),
    @("004641", "万家量化睿选灵活配置混合A", "5.24", "90.14", "1.36", "0.0713", 4),
    @("009206", "兴银丰运稳益回报混合C", "3.19", "38.30", "1.29", "0.0412", 8),
    @("002271", "招商安弘灵活配置混合", "0.55", "53.72", "5.02", "0.0276", 3),
    @("009205", "兴银丰运稳益回报混合A", "2.10", "38.30", "1.29", "0.0271", 8),
    @("016556", "万家量化睿选灵活配置混合C", "0.64", "90.14", "1.36", "0.0087", 4),
    @("005146", "兴银丰润灵活配置混合", "0.04", "92.81", "3.15", "0.0013", 9)
)

$r = 2
foreach ($row in $q3rows) {
    $q3.Cells.Item($r, 1).Value = $r - 2

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 2).Style = "Normal"

    $q3.Cells.Item($r, 3).Value = $row[1]

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 4).Style = "Normal"

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 5).Style = "Normal"

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 6).Style = "Normal"

    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 7).Style = "Normal"

    $q3.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add a 2022-Q3 row at the top of the
#    data (row 2) and push the existing quarters down by one row.
# ---------------------------------------------------------------------------
$ws = $zongji

# Make sure the whole index column (A2:A5) carries the bold/border/center
# look (some of these rows are brand-new cells with no inherited style).
$wsIndex = $ws.Range("A2:A5")
$wsIndex.Font.Bold = $true
$wsIndex.Borders.LineStyle = 1
$wsIndex.HorizontalAlignment = -4108
$wsIndex.VerticalAlignment = -4160

$summaryRows = @(
    @("2022-Q3", 10, 0.5600000000000001),
    @("2022-Q2", 21, 1.62),
    @("2022-Q1", 15, 0.8100000000000001),
    @("2021-Q4", 1, 0.03)
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
